$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Energy Equipment & Services(32)", 0.5423198209279574),
    @("Containers & Packaging(12)", 0.4271537997120528),
    @("Electric Utilities(28)", 0.4204824388844848),
    @("Marine(15)", 0.4187471659623295),
    @("Road & Rail(22)", 0.4163434739926311),
    @("Multi-Utilities(18)", 0.4101314264982653),
    @("Auto Components(21)", 0.4048105841666986),
    @("Chemicals(51)", 0.4001645687494276),
    @("Machinery(85)", 0.3919310531092248),
    @("Air Freight & Logistics(11)", 0.3843055748893925),
    @("Insurance(75)", 0.3663951961680489),
    @("Metals & Mining(89)", 0.354268143137205),
    @("Consumer Finance(15)", 0.3459232203685237),
    @("Oil, Gas & Consumable Fuels(122)", 0.338065111906787),
    @("Life Sciences Tools & Services(19)", 0.3230085395246081),
    @("Building Products(23)", 0.3174692589469055),
    @("Banks(246)", 0.2848110563497159),
    @("Trading Companies & Distributors(25)", 0.2783591086449008),
    @("Capital Markets(75)", 0.2770544365119251),
    @("Electrical Equipment(28)", 0.2747318155666205),
    @("Diversified Telecommunication Services(20)", 0.2696989538733799),
    @("Media(42)", 0.244822527871412),
    @("Semiconductors & Semiconductor Equipment(68)", 0.2392752732412778),
    @("Aerospace & Defense(37)", 0.2389711178610557),
    @("Construction & Engineering(20)", 0.2188648845285672),
    @("Commercial Services & Supplies(52)", 0.2169784159426315),
    @("Biotechnology(126)", 0.1995255792216599),
    @("Household Durables(39)", 0.1903255716778237),
    @("IT Services(52)", 0.1894672843098727),
    @("Health Care Providers & Services(46)", 0.1860657929027088),
    @("Pharmaceuticals(48)", 0.1803163295256444),
    @("Hotels, Restaurants & Leisure(50)", 0.1795038682429143),
    @("Thrifts & Mortgage Finance(47)", 0.1757322938026915),
    @("Software(66)", 0.1725878662950896),
    @("Health Care Equipment & Supplies(83)", 0.1621290761695298),
    @("Specialty Retail(58)", 0.1314788714608264),
    @("Professional Services(35)", 0.1221762510538849),
    @("Communications Equipment(45)", 0.1031403075966969)
)

$rowCount = $data.Length
for ($i = 0; $i -lt $rowCount; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
}

# Remove leftover rows that no longer have data (old sheet had up to row 41)
$lastRow = $rowCount + 2
$ws.Range("A" + $lastRow + ":B41").Clear()
